$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the old "create new account" failure-log row (row 3) into three
# separate rows, one per account type (FS / RS / PS). Insert two blank
# rows above the "valid user" row (currently row 4) to make room, which
# pushes it (and the notes row below it) down to rows 6 and 7.
$ws.Rows("4:5").Insert()

$ws.Range("B3").Value = "Logging in with invalid user, create new FS account"
$ws.Range("C3").Value = "N"
$ws.Range("D3").Value = "N/A"

$ws.Range("B4").Value = "Logging in with invalid user, create new RS account"
$ws.Range("C4").Value = "N"
$ws.Range("D4").Value = "N/A"

$ws.Range("B5").Value = "Logging in with invalid user, create new PS account"
$ws.Range("C5").Value = "N"
$ws.Range("D5").Value = "N/A"

$ws.Range("B6").Select()
